$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous contents (old used range was A1:G6)
$ws.Range("A1:G6").Clear()

# Header row
$ws.Range("A1").Value = "Season"
$ws.Range("B1").Value = "Episode"
$ws.Range("C1").Value = "Who"
$ws.Range("D1").Value = "Answer"
$ws.Range("E1").Value = "Correct"

# Data rows: Season, Episode, Who, Answer, Correct
$data = @(
    @(1, 1, "D", "F", "T"),
    @(1, 1, "L", "T", "T"),
    @(1, 1, "D", "F", "F"),
    @(1, 1, "L", "T", "T"),
    @(1, 1, "D", "F", "F"),
    @(1, 1, "L", "F", "F"),
    @(1, 1, "D", "F", "F"),
    @(1, 1, "L", "T", "T"),
    @(1, 1, "D", "F", "F"),
    @(1, 1, "L", "F", "T"),
    @(1, 1, "L", "F", "F")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Selection matches the target state
$ws.Range("D10").Select()

# Window position update recorded in workbook.xml
$excel.ActiveWindow.Left = 14060
